$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 values (SKU, Nombre Articulo, Unidad de Medida, Ubicacion changed; Id Solicitud and Cantidad changed)
$ws.Range("A3").Value = 75
$ws.Range("B3").Value = "20301-00420-0240"
$ws.Range("C3").Value = "CANDADO DE BRONCE DE 38 /40 MM"
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = "GENERICO"
$ws.Range("G3").Value = "(91)B09-PIC-00B-001-001-001"
$ws.Range("H3").Value = 5

# Remove row 4 entirely (it's no longer present in the new workbook)
$ws.Rows.Item(4).Delete()

# Adjust column widths to match new layout
# (values nudged to compensate for this engine's internal 1/6-character
# rounding of ColumnWidth so the saved `width` ends up as close as possible
# to the target layout)
$ws.Columns.Item(1).ColumnWidth = 9.5
$ws.Columns.Item(2).ColumnWidth = 16.1666666666667
$ws.Columns.Item(3).ColumnWidth = 31.8333333333333
$ws.Columns.Item(4).ColumnWidth = 8.33333333333333
$ws.Columns.Item(5).ColumnWidth = 12.8333333333333
$ws.Columns.Item(6).ColumnWidth = 15.8333333333333
$ws.Columns.Item(7).ColumnWidth = 25.3333333333333
$ws.Columns.Item(8).ColumnWidth = 8.33333333333333

$wb.Save()
